# Insert a new data row at row 980 (pushes existing rows 980-1049 down to 981-1050),
# mirroring the rest of the "Terminal La Palmera de La Serena - Coliflor" block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 980, shifting rows 980:1049 down to 981:1050.
$ws.Rows.Item(980).Insert()

# Populate the new row 980 with the new record.
$ws.Range("A980").Value2 = 8
$ws.Range("B980").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C980").Value2 = "Coquimbo"
$ws.Range("D980").Value2 = 45021
$ws.Range("D980").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E980").Value2 = 4
$ws.Range("F980").Value2 = 100112008
$ws.Range("G980").Value2 = "Coliflor"
$ws.Range("H980").Value2 = "Sin especificar"
$ws.Range("I980").Value2 = "Segunda"
$ws.Range("J980").Value2 = 1500
$ws.Range("K980").Value2 = 900
$ws.Range("L980").Value2 = 1000
$ws.Range("M980").Value2 = 950
$ws.Range("N980").Value2 = "$/unidad"
$ws.Range("O980").Value2 = "Provincia del Elquí"
$ws.Range("P980").Value2 = 950
$ws.Range("Q980").Value2 = 1
$ws.Range("R980").Value2 = "Hortaliza"
